# Insert a new weekly price record for Brócoli at
# Terminal Hortofrutícola Agro Chillán. This pushes the existing
# records that were in rows 169-207 down by one row (to 170-208)
# and fills in the brand-new row 169 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 169 (shifts rows 169:207 -> 170:208)
$ws.Rows.Item(169).Insert()

# Populate the new row 169 with the new record
$ws.Cells.Item(169, 1).Value = 7
$ws.Cells.Item(169, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value = "Ñuble"
$ws.Cells.Item(169, 4).Value = 44551
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = 100112023
$ws.Cells.Item(169, 7).Value = "Brócoli"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 2000
$ws.Cells.Item(169, 11).Value = 600
$ws.Cells.Item(169, 12).Value = 650
$ws.Cells.Item(169, 13).Value = 625
$ws.Cells.Item(169, 14).Value = "$/unidad"
$ws.Cells.Item(169, 15).Value = "Región del Maule"
$ws.Cells.Item(169, 16).Value = 625
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"
